$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at T (shifts the old "Faltas"/T column and its data to U)
$ws.Columns("T").Insert()

# 2) New class date header: 6-jun (serial 45083), same style/format as the other date headers
$ws.Range("T4").Value2 = 45083

# 3) Mark attendance ("F" = Falta/absent) for the new 6/jun class in the rows that have it
$rowsWithF = @(5,10,13,15,16,18,19,20,21,22,24,29,31,33,35,36,37)
foreach ($r in $rowsWithF) {
  $ws.Cells.Item($r, 20).Value2 = "F"
}

# 4) Update the "Faltas" totals column (now U) to include the new T column in its COUNTIF range
$ws.Range("U5:U37").Formula = '=COUNTIF(I5:T5,"F")'

# 5) Re-point the existing conditional formatting rules so they cover the new column layout
#    a) "F" highlighting across the attendance block now reaches column S
$cfF = $ws.Range("I5:R37").FormatConditions.Item(1)
$cfF.ModifyAppliesToRange($ws.Range("I5:S37"))

#    b) the two "lessThan" rules over I5:S37 stay on the same range (no change needed),
#       but make sure they still apply correctly
$cfLess = $ws.Range("I5:S37").FormatConditions

#    c) the "greaterThan 4" rule on the Faltas total now targets U5:U37 instead of T5:T37
$cfTotal = $ws.Range("T5:T37").FormatConditions.Item(1)
$cfTotal.ModifyAppliesToRange($ws.Range("U5:U37"))

# 6) Add the three conditional formatting rules for the new T column (copied visual style
#    from the equivalent rules used elsewhere in the sheet: red "F", green " ", green 0)
$newRuleF = $ws.Range("T5:T37").FormatConditions.Add(1, 3, '"F"')
$newRuleF.Font.Color = 6299648
$newRuleF.Interior.Color = 13551615

$newRuleBlank = $ws.Range("T5:T37").FormatConditions.Add(1, 3, '" "')
$newRuleBlank.Font.Color = 24832
$newRuleBlank.Interior.Color = 13561798

$newRuleZero = $ws.Range("T5:T37").FormatConditions.Add(1, 3, "0")
$newRuleZero.Font.Color = 24832
$newRuleZero.Interior.Color = 13561798

# 7) Selection cursor moved slightly (as saved by the author)
$ws.Range("V9").Select()
